# Update the staff list (Ho_va_ten / column C) with the new names
# reviewed with anh Manh, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newNames = @{
    2  = "Nguyễn Trung Trực"
    3  = "Đặng Anh Tuấn"
    4  = "Đặng Thu Thảo"
    5  = "Hoàng Nguyệt Anh"
    6  = "Công Tôn Sách"
    7  = "Trương Phi"
    8  = "Quan Vân Trường"
    9  = "Lưu Bị"
    10 = "Thủy Kính"
    11 = "Xích Thố"
    12 = "Hoàng Linh Mai"
    13 = "Tiểu Kiều"
}

foreach ($row in $newNames.Keys) {
    $ws.Cells.Item($row, 3).Value = $newNames[$row]
}

$ws.Range("C16").Select()
